# Results/ResultsTest.xlsx edit
# - Drop the unneeded D ("Dificultate detectie") and E ("Scor") columns.
# - Drop the discarded TEST*.jpg rows (TEST.jpg, TEST1-3.jpg, TEST5-8.jpg, TEST10.jpg),
#   keeping TEST4.jpg, TEST9.jpg and TEST11.jpg..TEST28.jpg.
# - Drop the trailing empty rows 144-152.
# - Add a running index (1..34) in column A next to every data row.
# - Give column B its own (wider) width.
# - Re-point the selection to E10 and recompute the average formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- remove the two trailing header columns (D, E) ----
$ws.Columns.Item(4).Clear()
$ws.Columns.Item(5).Clear()

# ---- remove the rows for the discarded test images (bottom-up so row
#      numbers of the rows still to be deleted stay valid) ----
$rowsToDelete = @(26,24,23,22,21,19,18,17,16)
foreach ($r in $rowsToDelete) {
  $ws.Rows.Item($r).Delete()
}

# ---- remove the trailing empty rows (144-152), bottom-up ----
for ($r = 152; $r -ge 144; $r--) {
  $ws.Rows.Item($r).Delete()
}

# ---- add the running index column (A2:A35 = 1..34) ----
for ($r = 2; $r -le 35; $r++) {
  $ws.Cells.Item($r, 1).Value = $r - 1
}

# ---- give column B its own width ----
$ws.Columns.Item(2).ColumnWidth = 16.25

# ---- selection / active cell ----
$ws.Range("E10").Select()
